$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2/15/2022"
$ws.Range("B8").Value = 0.05
$ws.Range("C8").Value = 0.15

$ws.Range("D8").Select()
